# Update "想去人数" (interested-count) figures on the 展览 and 全部类型 sheets
# to reflect the latest scrape (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value = 11509
$ws1.Range("F9").Value = 4343
$ws1.Range("F18").Value = 2212
$ws1.Range("F21").Value = 11296
$ws1.Range("F22").Value = 11206
$ws1.Range("F24").Value = 40

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 11509
$ws4.Range("F9").Value = 4343
$ws4.Range("F19").Value = 2212
$ws4.Range("F22").Value = 11296
$ws4.Range("F23").Value = 11206
$ws4.Range("F25").Value = 40
